# Arena page was renamed. Datetimepicker on department page was corrected.
#
# - "location-page.html" sheet is renamed to "arena-page.html".
# - The active/selected sheet moves from "competition-page.html"
#   (previously active) to the renamed "arena-page.html" sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("location-page.html")
$ws.Name = "arena-page.html"

# Make the (renamed) arena page the active tab, which also flips
# tabSelected off on the previously-active sheet (competition-page.html).
$ws.Activate()
